# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value  = 630
    $ws.Range("F5").Value  = 4591
    $ws.Range("F7").Value  = 130
    $ws.Range("F9").Value  = 3087
    $ws.Range("F12").Value = 250
    $ws.Range("F13").Value = 608

    if ($sheetName -eq "展览") {
        $ws.Range("F18").Value = 1765
        $ws.Range("F21").Value = 1571
        $ws.Range("F23").Value = 607
        $ws.Range("F28").Value = 93
        $ws.Range("F31").Value = 3678
        $ws.Range("F32").Value = 754
        $ws.Range("F34").Value = 451
        $ws.Range("F36").Value = 1760
    }
    elseif ($sheetName -eq "全部类型") {
        $ws.Range("F19").Value = 1765
        $ws.Range("F22").Value = 1571
        $ws.Range("F24").Value = 607
        $ws.Range("F29").Value = 93
        $ws.Range("F32").Value = 3678
        $ws.Range("F34").Value = 754
        $ws.Range("F36").Value = 451
        $ws.Range("F38").Value = 1760
    }
}
